# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 1016 (shifting existing rows 1016-1105
# down to 1017-1106) on the single "Fruta, Vega Modelo de Temuco -
# Mandarina" worksheet, then populate the newly inserted row with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 1016..1105 down by one row, creating a blank row 1016.
$ws.Rows(1016).Insert()

# Populate the newly inserted row 1016 with the new record.
$ws.Cells.Item(1016, 1).Value  = 10
$ws.Cells.Item(1016, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(1016, 3).Value  = "La Araucanía"
$ws.Cells.Item(1016, 4).Value  = 45223
$ws.Cells.Item(1016, 5).Value  = 9
$ws.Cells.Item(1016, 6).Value  = "Fruta"
$ws.Cells.Item(1016, 7).Value  = 100102
$ws.Cells.Item(1016, 8).Value  = "Cítricos"
$ws.Cells.Item(1016, 9).Value  = 100102004
$ws.Cells.Item(1016, 10).Value = "Mandarina"
$ws.Cells.Item(1016, 11).Value = "Murcott"
$ws.Cells.Item(1016, 12).Value = "Segunda"
$ws.Cells.Item(1016, 13).Value = 3
$ws.Cells.Item(1016, 14).Value = 200000
$ws.Cells.Item(1016, 15).Value = 200000
$ws.Cells.Item(1016, 16).Value = 200000
$ws.Cells.Item(1016, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(1016, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1016, 19).Value = 444
$ws.Cells.Item(1016, 20).Value = 450
